$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the data row (A3:H3) - values and hyperlink, keep styles/formatting
$ws.Range("A3:H3").ClearContents()

# Update selection to A3:H3 with active cell A3
$ws.Range("A3:H3").Select()
